$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet so the template placeholder drives the sheet name
$ws.Name = "{d.i18n.sheetName}"

# Localize header row labels (row 1) using i18n placeholders
$ws.Range("A1").Value = "{d.i18n.name}"
$ws.Range("B1").Value = "{d.i18n.address}"
$ws.Range("C1").Value = "{d.i18n.unitName}"
$ws.Range("D1").Value = "{d.i18n.phone}"
$ws.Range("E1").Value = "{d.i18n.email}"
